# Regenerate the "K" (strikeouts) column (column G) values for the
# cotton_jharel 2022 save_data sheet. The previous values were computed
# against the old "Strike#" metric; this writes the recalculated K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 3
    6  = 1
    7  = 3
    8  = 2
    9  = 1
    10 = 0
    11 = 0
    12 = 1
    13 = 4
    14 = 3
    15 = 1
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 2
    21 = 1
    22 = 1
    23 = 1
    24 = 0
    25 = 1
    26 = 1
    27 = 2
    28 = 1
    29 = 4
    30 = 1
    31 = 1
    32 = 1
    34 = 3
    35 = 3
    36 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
